$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated TPM-derived values to the data cells (rows 2-11).
# Values are quoted so the COM interop parses them as numeric literals
# (avoids a parser issue with bare scientific-notation numbers like 1E-05).
$ws.Range("I2").Value = "0.8757151715610434"
$ws.Range("J2").Value = "0.8757151715610434"
$ws.Range("M2").Value = "0.513343"
$ws.Range("N2").Value = "1.540029"
$ws.Range("O2").Value = "0.02896400434489499"
$ws.Range("P2").Value = "0.02896400434489499"
$ws.Range("Q2").Value = "0.1215515800263333"
$ws.Range("R2").Value = "1.093964220237"
$ws.Range("S2").Value = "0.02536421803398452"
$ws.Range("T2").Value = "0.02536421803398452"
$ws.Range("I3").Value = "0.8757151715610434"
$ws.Range("J3").Value = "0.8757151715610434"
$ws.Range("O3").Value = "0.8830650561604291"
$ws.Range("P3").Value = "0.8830650561604291"
$ws.Range("S3").Value = "0.7733134671550925"
$ws.Range("T3").Value = "0.7733134671550925"
$ws.Range("I4").Value = "0.8757151715610434"
$ws.Range("J4").Value = "0.8757151715610434"
$ws.Range("K4").Value = "2"
$ws.Range("L4").Value = "0.6666666666666666"
$ws.Range("M4").Value = "0.04863733333333333"
$ws.Range("N4").Value = "0.145912"
$ws.Range("O4").Value = "0.00274423131121058"
$ws.Range("P4").Value = "0.00274423131121058"
$ws.Range("Q4").Value = "0.01151655854844444"
$ws.Range("R4").Value = "0.103649026936"
$ws.Range("S4").Value = "0.00240316499349996"
$ws.Range("T4").Value = "0.00240316499349996"
$ws.Range("I5").Value = "0.8757151715610434"
$ws.Range("J5").Value = "0.8757151715610434"
$ws.Range("M5").Value = "1.499596"
$ws.Range("N5").Value = "4.498788"
$ws.Range("O5").Value = "0.08461068926543686"
$ws.Range("P5").Value = "0.08461068926543686"
$ws.Range("Q5").Value = "0.3550808391293334"
$ws.Range("R5").Value = "3.195727552164"
$ws.Range("S5").Value = "0.07409486426598016"
$ws.Range("T5").Value = "0.07409486426598016"
$ws.Range("I6").Value = "0.8757151715610434"
$ws.Range("J6").Value = "0.8757151715610434"
$ws.Range("K6").Value = "2"
$ws.Range("L6").Value = "0.6666666666666666"
$ws.Range("M6").Value = "0.010918"
$ws.Range("N6").Value = "0.032754"
$ws.Range("O6").Value = "0.0006160189180286153"
$ws.Range("P6").Value = "0.0006160189180286153"
$ws.Range("Q6").Value = "0.002585211351333333"
$ws.Range("R6").Value = "0.023266902162"
$ws.Range("S6").Value = "0.0005394571124862772"
$ws.Range("T6").Value = "0.0005394571124862772"
$ws.Range("G7").Value = "0.03360533333333333"
$ws.Range("H7").Value = "0.100816"
$ws.Range("I7").Value = "0.1242848284389566"
$ws.Range("J7").Value = "0.1242848284389566"
$ws.Range("M7").Value = "0.513343"
$ws.Range("N7").Value = "1.540029"
$ws.Range("O7").Value = "0.02896400434489499"
$ws.Range("P7").Value = "0.02896400434489499"
$ws.Range("Q7").Value = "0.01725106262933333"
$ws.Range("R7").Value = "0.155259563664"
$ws.Range("S7").Value = "0.003599786310910468"
$ws.Range("T7").Value = "0.003599786310910468"
$ws.Range("G8").Value = "0.03360533333333333"
$ws.Range("H8").Value = "0.100816"
$ws.Range("I8").Value = "0.1242848284389566"
$ws.Range("J8").Value = "0.1242848284389566"
$ws.Range("O8").Value = "0.8830650561604291"
$ws.Range("P8").Value = "0.8830650561604291"
$ws.Range("Q8").Value = "0.5259566463324444"
$ws.Range("R8").Value = "4.733609816992"
$ws.Range("S8").Value = "0.1097515890053365"
$ws.Range("T8").Value = "0.1097515890053365"
$ws.Range("G9").Value = "0.03360533333333333"
$ws.Range("H9").Value = "0.100816"
$ws.Range("I9").Value = "0.1242848284389566"
$ws.Range("J9").Value = "0.1242848284389566"
$ws.Range("K9").Value = "2"
$ws.Range("L9").Value = "0.6666666666666666"
$ws.Range("M9").Value = "0.04863733333333333"
$ws.Range("N9").Value = "0.145912"
$ws.Range("O9").Value = "0.00274423131121058"
$ws.Range("P9").Value = "0.00274423131121058"
$ws.Range("Q9").Value = "0.001634473799111111"
$ws.Range("R9").Value = "0.014710264192"
$ws.Range("S9").Value = "0.0003410663177106198"
$ws.Range("T9").Value = "0.0003410663177106199"
$ws.Range("G10").Value = "0.03360533333333333"
$ws.Range("H10").Value = "0.100816"
$ws.Range("I10").Value = "0.1242848284389566"
$ws.Range("J10").Value = "0.1242848284389566"
$ws.Range("M10").Value = "1.499596"
$ws.Range("N10").Value = "4.498788"
$ws.Range("O10").Value = "0.08461068926543686"
$ws.Range("P10").Value = "0.08461068926543686"
$ws.Range("Q10").Value = "0.05039442344533334"
$ws.Range("R10").Value = "0.453549811008"
$ws.Range("S10").Value = "0.01051582499945669"
$ws.Range("T10").Value = "0.01051582499945669"
$ws.Range("G11").Value = "0.03360533333333333"
$ws.Range("H11").Value = "0.100816"
$ws.Range("I11").Value = "0.1242848284389566"
$ws.Range("J11").Value = "0.1242848284389566"
$ws.Range("K11").Value = "2"
$ws.Range("L11").Value = "0.6666666666666666"
$ws.Range("M11").Value = "0.010918"
$ws.Range("N11").Value = "0.032754"
$ws.Range("O11").Value = "0.0006160189180286153"
$ws.Range("P11").Value = "0.0006160189180286153"
$ws.Range("Q11").Value = "0.0003669030293333333"
$ws.Range("R11").Value = "0.003302127264"
$ws.Range("S11").Value = "7.656180554233813E-05"
$ws.Range("T11").Value = "7.656180554233814E-05"

# A couple of very small magnitude values get an auto-applied scientific
# number format from the interop; reset those specific cells back to the
# default "Normal" style so only the values (not formatting) change.
$ws.Range("S11").Style = "Normal"
$ws.Range("T11").Style = "Normal"
